# Workbook was edited after adjusting the ML "Label" (Buy signal) column.
# The sheet with the trading calc data is stored internally as "Sheet2"
# (its physical part is worksheets/sheet1.xml) and is the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column C ("Label" / buy-signal flag), keyed by row number.
# Only rows whose flag actually changed are listed; Excel will recompute the
# dependent D/E/H/I/J formula columns automatically.
$newLabels = @{
    2  = 1
    3  = 0
    4  = 0
    5  = 1
    12 = 1
    14 = 1
    15 = 1
    16 = 0
    19 = 0
    20 = 1
    21 = 0
    22 = 1
    24 = 0
    25 = 1
    26 = 0
    27 = 1
    28 = 0
    29 = 1
    30 = 1
    31 = 0
    34 = 1
    35 = 1
    39 = 1
    40 = 1
    41 = 1
    43 = 0
    45 = 0
    46 = 0
    47 = 1
    48 = 1
}

foreach ($row in $newLabels.Keys) {
    $ws.Cells.Item($row, 3).Value = $newLabels[$row]
}

# Touch the final formula cells so the recalculation engine evaluates the
# whole dependency chain (D/E columns down to row 51, and H2/H3/I2/I3/J3).
$null = $ws.Range("E51").Value2
$null = $ws.Range("H2").Value2
$null = $ws.Range("H3").Value2
$null = $ws.Range("I2").Value2
$null = $ws.Range("I3").Value2
$null = $ws.Range("J3").Value2
